$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4 is an untouched, default-styled ("General" number format, no explicit style)
# cell used purely as a format donor: after forcing a column-D cell to Text
# (via NumberFormat "@") so Excel stores e.g. "582.82" as a string instead of
# silently parsing it into a Number, we paste D4's format back on top so the
# cell keeps its original (default) style/number format.
$donor = $ws.Range("D4")

$ws.Range("D2").Value = '67.132.34'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '2.470.37'
$ws.Range("E3").Value = '  -1.84%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.82'
$donor.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.31'
$donor.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("D9").Value = '2.470.07'
$ws.Range("E9").Value = '  -1.76%  '
$ws.Range("E10").Value = '  -2.34%  '
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.53'
$donor.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = '  -3.29%  '
$ws.Range("D16").Value = '67.058.08'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000169'
$donor.Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = '  -4.08%  '
$ws.Range("D18").Value = '2.450.51'
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.16'
$donor.Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = '  -5.58%  '
$ws.Range("E20").Value = '  -4.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.24'
$donor.Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = '  -3.56%  '
$ws.Range("E22").Value = '  -2.79%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.18'
$donor.Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E25").Value = '  -7.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.79'
$donor.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = '  -6.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.19'
$donor.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = '  -8.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$donor.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").Value = '2.593.37'
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0906'
$ws.Range("E30").Value = '  -5.46%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '515.04'
$donor.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = '  -3.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.74'
$donor.Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = '  -7.13%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$donor.Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = '  -5.67%  '
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.78'
$donor.Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("E34").Value = '  -5.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$donor.Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -6.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.56'
$donor.Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.39'
$donor.Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("E40").Value = '  -5.41%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.67'
$donor.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = '  -6.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.327'
$donor.Copy()
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("E43").Value = '  -6.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.80'
$donor.Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = '  -6.33%  '
$ws.Range("E45").Value = '  -5.12%  '
$ws.Range("E46").Value = '  -2.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.93'
$donor.Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("E47").Value = '  -3.57%  '
$ws.Range("E48").Value = '  -6.35%  '
$ws.Range("E49").Value = '  -6.68%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0252'
$ws.Range("E50").Value = '  -13.26%  '
$ws.Range("B51").Value = 'Optimism'
$ws.Range("C51").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.59'
$donor.Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = '  -7.19%  '

$excel.CutCopyMode = 0
